{"js": "// Increment 2 Progress Report \u2014 minor correction pass.\n//\n// 1) Insert a new sentence fragment into William Hudmon's \"Progress\n//    Report\" contribution bullet (the only genuine content change).\n// 2) Normalize the document body: this drops the stray w:proofErr\n//    (spell/grammar \"squiggle\") markers and coalesces runs that are\n//    split purely because of those markers or other incidental\n//    authoring splits \u2014 mirroring what Word itself does when the\n//    content is re-serialized. The visible text is unaffected except\n//    for the sentence inserted in step 1.\n\nconst body = context.document.body;\n\n// --- Step 1: add the missing contribution note -------------------------\nconst oldSentence =\n  \"Progress Report: Gave input/suggestions and discussion for the sections of the document relevant to Increment 2.\";\nconst newSentence =\n  \"Progress Report: Gave input/suggestions and discussion for the sections of the document relevant to Increment 2; Wrote the entirety of the Stakeholder Communication letter.\";\n\nconst target = body.search(oldSentence, { matchCase: true });\ntarget.load(\"text\");\nawait context.sync();\n\nif (target.items.length !== 1) {\n  throw new Error(\n    \"Expected exactly one match for the William Hudmon progress-report sentence, found \" +\n      target.items.length\n  );\n}\ntarget.items[0].insertText(newSentence, \"Replace\");\nawait context.sync();\n\n// --- Step 2: normalize runs / strip proofErr markers --------------------\nconst wholeRange = body.getRange(\"Whole\");\nconst ooxml = body.getOoxml();\nawait context.sync();\n\nwholeRange.insertOoxml(ooxml.value, \"Replace\");\nawait context.sync();\n", "ps1": "# Increment 2 Progress Report -- minor correction pass.\n#\n# 1) Insert a new sentence fragment into William Hudmon's \"Progress\n#    Report\" contribution bullet (the only genuine content change).\n# 2) Normalize the document body: this drops the stray proofErr\n#    (spell/grammar \"squiggle\") markers and coalesces runs that are\n#    split purely because of those markers or other incidental\n#    authoring splits -- mirroring what Word itself does when the\n#    content is re-serialized through WordOpenXML. The visible text is\n#    unaffected except for the sentence inserted in step 1.\n\n$d = $word.ActiveDocument\n\n# --- Step 1: add the missing contribution note --------------------------\n$oldSentence = \"Progress Report: Gave input/suggestions and discussion for the sections of the document relevant to Increment 2.\"\n$newSentence = \"Progress Report: Gave input/suggestions and discussion for the sections of the document relevant to Increment 2; Wrote the entirety of the Stakeholder Communication letter.\"\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Text = $oldSentence\n$found = $find.Execute()\nif (-not $found) {\n    throw \"Could not find the William Hudmon progress-report sentence to update.\"\n}\n$find.Parent.Text = $newSentence\n\n# --- Step 2: normalize runs / strip proofErr markers ---------------------\n$xml = $d.Content.WordOpenXML\n$whole = $d.Range()\n$whole.InsertXML($xml)\n"}
